$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Inscritos 12 -> 13
$ws.Range("E4").Value = 13

# Row 15: Inscritos 53 -> 54
$ws.Range("E15").Value = 54

# Row 16: Inscritos 205 -> 206
$ws.Range("E16").Value = 206

# Row 17: Inscritos 11 -> 13, Pagos 2 -> 3, Inscricoes homologadas 2 -> 3
$ws.Range("E17").Value = 13
$ws.Range("F17").Value = 3
$ws.Range("H17").Value = 3

# Row 18: Inscritos 46 -> 47
$ws.Range("E18").Value = 47
